$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 115.4136146666667
$ws.Range("H2").Value = 346.240844
$ws.Range("I2").Value = 0.2619217538490851
$ws.Range("J2").Value = 0.2619217538490851
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.229822
$ws.Range("N2").Value = 0.689466
$ws.Range("O2").Value = 0.09226175421862418
$ws.Range("P2").Value = 0.09226175421862419
$ws.Range("Q2").Value = 26.52458774992267
$ws.Range("R2").Value = 238.721289749304
$ws.Range("S2").Value = 0.02416536047813527
$ws.Range("T2").Value = 0.02416536047813527
$ws.Range("G3").Value = 115.4136146666667
$ws.Range("H3").Value = 346.240844
$ws.Range("I3").Value = 0.2619217538490851
$ws.Range("J3").Value = 0.2619217538490851
$ws.Range("O3").Value = 0.4364142651333466
$ws.Range("P3").Value = 0.4364142651333466
$ws.Range("Q3").Value = 125.4659481481089
$ws.Range("R3").Value = 1129.19353333298
$ws.Range("S3").Value = 0.1143063897284858
$ws.Range("T3").Value = 0.1143063897284858
$ws.Range("G4").Value = 115.4136146666667
$ws.Range("H4").Value = 346.240844
$ws.Range("I4").Value = 0.2619217538490851
$ws.Range("J4").Value = 0.2619217538490851
$ws.Range("M4").Value = 1.174057666666666
$ws.Range("O4").Value = 0.4713239806480292
$ws.Range("P4").Value = 0.4713239806480293
$ws.Range("Q4").Value = 135.5022391371124
$ws.Range("R4").Value = 1219.520152234012
$ws.Range("S4").Value = 0.123450003642464
$ws.Range("T4").Value = 0.1234500036424641
$ws.Range("I5").Value = 0.6414314537852458
$ws.Range("J5").Value = 0.6414314537852458
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.229822
$ws.Range("N5").Value = 0.689466
$ws.Range("O5").Value = 0.09226175421862418
$ws.Range("P5").Value = 0.09226175421862419
$ws.Range("Q5").Value = 64.95720432328133
$ws.Range("R5").Value = 584.614838909532
$ws.Range("S5").Value = 0.05917959113722914
$ws.Range("T5").Value = 0.05917959113722915
$ws.Range("I6").Value = 0.6414314537852458
$ws.Range("J6").Value = 0.6414314537852458
$ws.Range("O6").Value = 0.4364142651333466
$ws.Range("P6").Value = 0.4364142651333466
$ws.Range("R6").Value = 2765.330634232089
$ws.Range("S6").Value = 0.2799298365371022
$ws.Range("T6").Value = 0.2799298365371022
$ws.Range("I7").Value = 0.6414314537852458
$ws.Range("J7").Value = 0.6414314537852458
$ws.Range("M7").Value = 1.174057666666666
$ws.Range("O7").Value = 0.4713239806480292
$ws.Range("P7").Value = 0.4713239806480293
$ws.Range("Q7").Value = 331.8372642348495
$ws.Range("R7").Value = 2986.535378113645
$ws.Range("S7").Value = 0.3023220261109144
$ws.Range("T7").Value = 0.3023220261109145
$ws.Range("I8").Value = 0.09664679236566912
$ws.Range("J8").Value = 0.09664679236566913
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.229822
$ws.Range("N8").Value = 0.689466
$ws.Range("O8").Value = 0.09226175421862418
$ws.Range("P8").Value = 0.09226175421862419
$ws.Range("Q8").Value = 9.787336436089999
$ws.Range("R8").Value = 88.08602792481
$ws.Range("S8").Value = 0.008916802603259767
$ws.Range("T8").Value = 0.00891680260325977
$ws.Range("I9").Value = 0.09664679236566912
$ws.Range("J9").Value = 0.09664679236566913
$ws.Range("O9").Value = 0.4364142651333466
$ws.Range("P9").Value = 0.4364142651333466
$ws.Range("Q9").Value = 46.29581644684166
$ws.Range("R9").Value = 416.6623480215749
$ws.Range("S9").Value = 0.04217803886775862
$ws.Range("T9").Value = 0.04217803886775862
$ws.Range("I10").Value = 0.09664679236566912
$ws.Range("J10").Value = 0.09664679236566913
$ws.Range("M10").Value = 1.174057666666666
$ws.Range("O10").Value = 0.4713239806480292
$ws.Range("P10").Value = 0.4713239806480293
$ws.Range("Q10").Value = 49.99911835697832
$ws.Range("R10").Value = 449.9920652128049
$ws.Range("S10").Value = 0.04555195089465073
$ws.Range("T10").Value = 0.04555195089465074
